$wb = $excel.ActiveWorkbook

# Overview sheet: row 4 corresponds to file 6766224b-3f0f-458e-9f7a-bd1b75c51507.md
# "Latest HO Xliff Generate Date" (column G) gets a fresh handoff-report timestamp.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-10-17 14:01:20"

# zh-cn sheet: same row's "Latest Handoff Datetime" (column H) updated.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-10-17 14:00:58"

# de-de sheet: same row's "Latest Handoff Datetime" (column H) updated.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H4").Value = "2016-10-17 14:01:20"
